$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 1346
$ws.Range("I2").Value = 3507
$ws.Range("J2").Value = 14564
$ws.Range("K2").Value = 80
$ws.Range("L2").Value = 4061
$ws.Range("M2").Value = 240
$ws.Range("N2").Value = 2369
$ws.Range("O2").Value = 5
$ws.Range("P2").Value = 61
$ws.Range("Q2").Value = 31
$ws.Range("R2").Value = 203
$ws.Range("S2").Value = 1540
$ws.Range("T2").Value = 2550
$ws.Range("U2").Value = 185
$ws.Range("V2").Value = 22453
$ws.Range("X2").Value = 22370
$ws.Range("Y2").Value = 34
$ws.Range("Z2").Value = 337
$ws.Range("AA2").Value = 151
